# Update "想去人数" (number of people interested) counts for a couple of
# entries that are duplicated across the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 495
    $ws.Range("F8").Value = 1981
    $ws.Range("F9").Value = 4077
    $ws.Range("F10").Value = 95
}
